$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# ---------------------------------------------------------------------------
# 1. Insert a new "KVM RAM" sub-header row above the old row 89 (KVM VCIN RAM)
#    This pushes everything from row 89 onward down by one row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(89).Insert()
$ws.Range("A85").Copy()
$ws.Range("A89").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A89:B89").Merge()
$ws.Range("B89").ClearFormats()
$ws.Range("A89").Value = "KVM RAM"

# ---------------------------------------------------------------------------
# 2. Insert a new "KVM CPU" sub-header row above what is now row 99
#    (originally row 98, "KVM NUH CPU cores", before the first insert).
# ---------------------------------------------------------------------------
$ws.Rows.Item(99).Insert()
$ws.Range("A85").Copy()
$ws.Range("A99").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A99:B99").Merge()
$ws.Range("B99").ClearFormats()
$ws.Range("A99").Value = "KVM CPU"

# ---------------------------------------------------------------------------
# 3. Rename the labels in column A (rows now at their final positions).
# ---------------------------------------------------------------------------
$ws.Range("A85").Value = "KVM and VCenter RAM"
$ws.Range("A86").Value = "KVM VSD RAM"
$ws.Range("A87").Value = "KVM VSC RAM"
$ws.Range("A88").Value = "KVM VSTAT RAM"

$ws.Range("A94").Value = "KVM and VCenter CPU"
$ws.Range("A95").Value = "KVM VSD CPU cores"
$ws.Range("A96").Value = "KVM VSC CPU cores"
$ws.Range("A97").Value = "KVM VSTAT CPU cores"
$ws.Range("A98").Value = "KVM VNSUTIL CPU cores"

# ---------------------------------------------------------------------------
# 4. Update comment text to match new wording.
# ---------------------------------------------------------------------------
$ws.Range("A86").Comment.Text("For KVM and VCenter deployments: amount of VSD RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]")
$ws.Range("A87").Comment.Text("For KVM and VCenter deployments: amount of VSC RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 4]")
$ws.Range("A88").Comment.Text("For KVM and VCenter deployments: amount of VSTAT RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 16]")

$ws.Range("A95").Comment.Text("For KVM and VCenter deployments: number of CPU's for VSD. [default: 6]")
$ws.Range("A96").Comment.Text("For KVM and VCenter deployments: number of CPU's for VSC. [default: 6]")
$ws.Range("A97").Comment.Text("For KVM and VCenter deployments: number of CPU's for VSTAT. [default: 6]")
$ws.Range("A98").Comment.Text("For KVM and VCenter deployments: number of CPU's for VNSUTIL. [default: 2]")
